$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E). This shifts all columns
# to the right of it (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) one position
# to the left.
$ws.Range("E1").EntireColumn.Delete()
